# ArmLab_P2.docx revision: replace the "Arm Lab P2" discussion
# paragraphs with the finished answers for 2.1 / 2.2 / 2.3, and drop
# the whole "Figure 1 config space" section (heading, caption, image,
# joint-limit text) that used to follow the page break.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Drop everything from the "2.1" heading paragraph (right after
#    the page-break paragraph) through the very end of the document
#    (figure heading, caption, picture, joint-limit paragraphs and
#    the trailing blank paragraph).
# ---------------------------------------------------------------
$tailStart = $d.Paragraphs.Item(7).Range.Start
$tailEnd   = $d.Content.End
$d.Range($tailStart, $tailEnd).Delete()

# ---------------------------------------------------------------
# 2) Drop the third bullet ("Sturdiness and precision ...") entirely.
# ---------------------------------------------------------------
$d.Paragraphs.Item(5).Range.Delete()

# ---------------------------------------------------------------
# 3) Paragraph 3 ("I would place the motor ...") -> answer to 2.1
# ---------------------------------------------------------------
$origP3 = "I would place the motor so that it rotates near each joint. The first motor could be placed anywhere. The second one would be placed along the 3.75 inch link. The motor would be attached to each joint with gears, so we can adjust the gear ratio. In addition, if the third dimension is considered, positioning the motors above the arm will prevent the motors from hitting the obstacles."
$d.Content.Find.Execute($origP3, $true, $false, $false, $false, $false, $true, 1, $false, "2.1 The viable points will be represented by grid", 2) | Out-Null

$p3 = $d.Paragraphs.Item(3)
$p3.Range.ListFormat.RemoveNumbers()
$p3.Style = "Normal"

$p3 = $d.Paragraphs.Item(3)
$p3EndIns = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$p3EndIns.InsertAfter(", because it is relatively easy to compute intermediate waypoints.")

# ---------------------------------------------------------------
# 4) Paragraph 4 ("A high gear ratio ...") -> answer to 2.2
# ---------------------------------------------------------------
$origP4 = "A high gear ratio is better than a low one, because speed is not as important as the torque in order to control the arms with accuracy. We would want a higher gear ratio at the first joint because more of the arm is cantilevered off of that point. The second joint should also have a high gear ratio, but it can be less than the first joint. "
$d.Content.Find.Execute($origP4, $true, $false, $false, $false, $false, $true, 1, $false, "2.2 The robot will reach the goals using Wavefront Planning. The biggest benefit to using Wavefront for this lab is because the configuration created for part 1 made it easy to define the world in part2 by converting the configuration space to a binary grid. We defined ones for obstacles and zeros for reachable space, then performed the Wavefront Planner. Besides, it was useful to use this method since we used Wavefront planning for the Motion Planning lab; it was the most convenient method to update the strategy to fit the arm-lab.", 2) | Out-Null

$p4 = $d.Paragraphs.Item(4)
$p4.Range.ListFormat.RemoveNumbers()
$p4.Style = "Normal"

# Put the (hidden) _GoBack bookmark back where the last bit of typing
# happened, right before "the strategy to fit the arm-lab."
$p4 = $d.Paragraphs.Item(4)
$p4Text = $p4.Range.Text
$tailPhrase = "the strategy to fit the arm-lab."
$bmOffset = $p4Text.IndexOf($tailPhrase)
if ($bmOffset -ge 0) {
    $bmPos = $p4.Range.Start + $bmOffset
    $d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos)) | Out-Null
}

# ---------------------------------------------------------------
# 5) The page-break paragraph becomes the answer to 2.3
# ---------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5End = $d.Range($p5.Range.End - 1, $p5.Range.End - 1)
$p5End.InsertAfter("2.3 Using two PID controls with a feed forward; one PID constant for each link.")
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Start
$breakRange = $d.Range($p5.Range.Start, $p5.Range.Start + 1)
$breakRange.Delete()
